$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title (A1) changes to reflect the "For Test" folder name
$ws.Range("A1").Value = "Benchmark For Test/BUFF Vulnerability Factor"

# Column F header changes from "Gate Inputs List" to "Gate Delay"
$ws.Range("F2").Value = "Gate Delay"

# Column F now holds the computed gate delay (numeric) instead of the
# textual list of input-gate ids
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("F5").Value = 19.8
$ws.Range("F6").Value = 29.87
$ws.Range("F7").Value = 39.94
$ws.Range("F8").Value = 50.01
